$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "model_2_7_24"
$ws.Cells.Item(2, 2).Value = 0.9735302659487451
$ws.Cells.Item(2, 3).Value = 0.914985832995417
$ws.Cells.Item(2, 4).Value = 0.9162931109510041
$ws.Cells.Item(2, 5).Value = 0.917003051648597
$ws.Cells.Item(2, 6).Value = 3.192257807431467
$ws.Cells.Item(2, 7).Value = 11.8566552488194
$ws.Cells.Item(2, 8).Value = 8.950469229650569
$ws.Cells.Item(2, 9).Value = 10.48903930983894

# Row 3
$ws.Cells.Item(3, 1).Value = "model_2_7_23"
$ws.Cells.Item(3, 2).Value = 0.9737900346362308
$ws.Cells.Item(3, 3).Value = 0.9158012232914017
$ws.Cells.Item(3, 4).Value = 0.9172715962389163
$ws.Cells.Item(3, 5).Value = 0.9178690545348773
$ws.Cells.Item(3, 6).Value = 3.160929626379597
$ws.Cells.Item(3, 7).Value = 11.7429353598484
$ws.Cells.Item(3, 8).Value = 8.845843403023586
$ws.Cells.Item(3, 9).Value = 10.37959506523648

# Row 4
$ws.Cells.Item(4, 1).Value = "model_2_7_22"
$ws.Cells.Item(4, 2).Value = 0.974074016374483
$ws.Cells.Item(4, 3).Value = 0.9166983622725982
$ws.Cells.Item(4, 4).Value = 0.9183506087303057
$ws.Cells.Item(4, 5).Value = 0.918822777078437
$ws.Cells.Item(4, 6).Value = 3.126681344196325
$ws.Cells.Item(4, 7).Value = 11.61781424197931
$ws.Cells.Item(4, 8).Value = 8.730468572919262
$ws.Cells.Item(4, 9).Value = 10.25906493191492

# Row 5
$ws.Cells.Item(5, 1).Value = "model_2_7_21"
$ws.Cells.Item(5, 2).Value = 0.9743839783856015
$ws.Cells.Item(5, 3).Value = 0.9176846213117951
$ws.Cells.Item(5, 4).Value = 0.9195397913747663
$ws.Cells.Item(5, 5).Value = 0.9198724732533931
$ws.Cells.Item(5, 6).Value = 3.089299833370254
$ws.Cells.Item(5, 7).Value = 11.48026383331436
$ws.Cells.Item(5, 8).Value = 8.603313654266762
$ws.Cells.Item(5, 9).Value = 10.12640578406428

# Row 6
$ws.Cells.Item(6, 1).Value = "model_2_7_20"
$ws.Cells.Item(6, 2).Value = 0.9747214230588095
$ws.Cells.Item(6, 3).Value = 0.91876751271101
$ws.Cells.Item(6, 4).Value = 0.9208484586917673
$ws.Cells.Item(6, 5).Value = 0.9210262342896733
$ws.Cells.Item(6, 6).Value = 3.048603905313759
$ws.Cells.Item(6, 7).Value = 11.32923641700491
$ws.Cells.Item(6, 8).Value = 8.463382679818404
$ws.Cells.Item(6, 9).Value = 9.980595063259619

# Row 7
$ws.Cells.Item(7, 1).Value = "model_2_7_19"
$ws.Cells.Item(7, 2).Value = 0.9750878293658223
$ws.Cells.Item(7, 3).Value = 0.919954769153413
$ws.Cells.Item(7, 4).Value = 0.9222875260295545
$ws.Cells.Item(7, 5).Value = 0.9222928366496776
$ws.Cells.Item(7, 6).Value = 3.004415195597649
$ws.Cells.Item(7, 7).Value = 11.16365354034444
$ws.Cells.Item(7, 8).Value = 8.309508511603642
$ws.Cells.Item(7, 9).Value = 9.820523612345907

# Row 8
$ws.Cells.Item(8, 1).Value = "model_2_7_18"
$ws.Cells.Item(8, 2).Value = 0.9754843779338311
$ws.Cells.Item(8, 3).Value = 0.9212540089107064
$ws.Cells.Item(8, 4).Value = 0.9238678080168059
$ws.Cells.Item(8, 5).Value = 0.9236811174615421
$ws.Cells.Item(8, 6).Value = 2.956591320231137
$ws.Cells.Item(8, 7).Value = 10.98245270722968
$ws.Cells.Item(8, 8).Value = 8.140534781224224
$ws.Cells.Item(8, 9).Value = 9.645074607316372

# Row 9
$ws.Cells.Item(9, 1).Value = "model_2_7_17"
$ws.Cells.Item(9, 2).Value = 0.9759119788504422
$ws.Cells.Item(9, 3).Value = 0.9226733355454517
$ws.Cells.Item(9, 4).Value = 0.9256005226614985
$ws.Cells.Item(9, 5).Value = 0.9252002453371605
$ws.Cells.Item(9, 6).Value = 2.905022522377946
$ws.Cells.Item(9, 7).Value = 10.78450373958602
$ws.Cells.Item(9, 8).Value = 7.955261988419699
$ws.Cells.Item(9, 9).Value = 9.453089331706364

# Row 10
$ws.Cells.Item(10, 1).Value = "model_2_7_16"
$ws.Cells.Item(10, 2).Value = 0.9763710441062221
$ws.Cells.Item(10, 3).Value = 0.9242202155070322
$ws.Cells.Item(10, 4).Value = 0.9274976783897
$ws.Cells.Item(10, 5).Value = 0.9268593150869346
$ws.Cells.Item(10, 6).Value = 2.849659115853114
$ws.Cells.Item(10, 7).Value = 10.56876531548575
$ws.Cells.Item(10, 8).Value = 7.752406116435451
$ws.Cells.Item(10, 9).Value = 9.243418395981514

# Row 11
$ws.Cells.Item(11, 1).Value = "model_2_7_15"
$ws.Cells.Item(11, 2).Value = 0.976861341939983
$ws.Cells.Item(11, 3).Value = 0.9259020146976372
$ws.Cells.Item(11, 4).Value = 0.9295704156992129
$ws.Cells.Item(11, 5).Value = 0.9286671901854622
$ws.Cells.Item(11, 6).Value = 2.790529051124883
$ws.Cells.Item(11, 7).Value = 10.33421013599816
$ws.Cells.Item(11, 8).Value = 7.530775952888392
$ws.Cells.Item(11, 9).Value = 9.014941646505786

# Row 12
$ws.Cells.Item(12, 1).Value = "model_2_7_14"
$ws.Cells.Item(12, 2).Value = 0.9773817156262434
$ws.Cells.Item(12, 3).Value = 0.9277252447475274
$ws.Cells.Item(12, 4).Value = 0.9318302967304162
$ws.Cells.Item(12, 5).Value = 0.930632173741483
$ws.Cells.Item(12, 6).Value = 2.727771829630703
$ws.Cells.Item(12, 7).Value = 10.07993004477914
$ws.Cells.Item(12, 8).Value = 7.289135200708298
$ws.Cells.Item(12, 9).Value = 8.766609747903614

# Row 13
$ws.Cells.Item(13, 1).Value = "model_2_7_13"
$ws.Cells.Item(13, 2).Value = 0.9779298656955852
$ws.Cells.Item(13, 3).Value = 0.929694775021601
$ws.Cells.Item(13, 4).Value = 0.9342887706050403
$ws.Cells.Item(13, 5).Value = 0.9327615953902633
$ws.Cells.Item(13, 6).Value = 2.661664768067023
$ws.Cells.Item(13, 7).Value = 9.805245927006794
$ws.Cells.Item(13, 8).Value = 7.026259647492553
$ws.Cells.Item(13, 9).Value = 8.497496391027997

# Row 14
$ws.Cells.Item(14, 1).Value = "model_2_7_12"
$ws.Cells.Item(14, 2).Value = 0.9785019361806913
$ws.Cells.Item(14, 3).Value = 0.9318142698889389
$ws.Cells.Item(14, 4).Value = 0.9369548184507606
$ws.Cells.Item(14, 5).Value = 0.9350615111009482
$ws.Cells.Item(14, 6).Value = 2.592672897240337
$ws.Cells.Item(14, 7).Value = 9.509646724790155
$ws.Cells.Item(14, 8).Value = 6.741188974349662
$ws.Cells.Item(14, 9).Value = 8.206836230891129

# Row 15
$ws.Cells.Item(15, 1).Value = "model_2_7_11"
$ws.Cells.Item(15, 2).Value = 0.9790920182241202
$ws.Cells.Item(15, 3).Value = 0.9340846740509849
$ws.Cells.Item(15, 4).Value = 0.9398371794182344
$ws.Cells.Item(15, 5).Value = 0.9375355778534525
$ws.Cells.Item(15, 6).Value = 2.521508827117319
$ws.Cells.Item(15, 7).Value = 9.193000683626053
$ws.Cells.Item(15, 8).Value = 6.432988736099023
$ws.Cells.Item(15, 9).Value = 7.894167103439476

# Row 16
$ws.Cells.Item(16, 1).Value = "model_2_7_10"
$ws.Cells.Item(16, 2).Value = 0.9796915974074193
$ws.Cells.Item(16, 3).Value = 0.9365036492689131
$ws.Cells.Item(16, 4).Value = 0.9429413353882211
$ws.Cells.Item(16, 5).Value = 0.9401846735866956
$ws.Cells.Item(16, 6).Value = 2.449199399098375
$ws.Cells.Item(16, 7).Value = 8.855633910240332
$ws.Cells.Item(16, 8).Value = 6.101072775428925
$ws.Cells.Item(16, 9).Value = 7.55937805596911

# Row 17
$ws.Cells.Item(17, 1).Value = "model_2_7_0"
$ws.Cells.Item(17, 2).Value = 0.9802444861292557
$ws.Cells.Item(17, 3).Value = 0.9624396333084851
$ws.Cells.Item(17, 4).Value = 0.9821355742006392
$ws.Cells.Item(17, 5).Value = 0.970942878233279
$ws.Cells.Item(17, 6).Value = 2.382520854632997
$ws.Cells.Item(17, 7).Value = 5.238424777561168
$ws.Cells.Item(17, 8).Value = 1.910177229606076
$ws.Cells.Item(17, 9).Value = 3.67219877954419

# Row 18
$ws.Cells.Item(18, 1).Value = "model_2_7_9"
$ws.Cells.Item(18, 2).Value = 0.9802886922022972
$ws.Cells.Item(18, 3).Value = 0.9390634927000943
$ws.Cells.Item(18, 4).Value = 0.9462688317936079
$ws.Cells.Item(18, 5).Value = 0.9430051322708402
$ws.Cells.Item(18, 6).Value = 2.377189589062696
$ws.Cells.Item(18, 7).Value = 8.498620695574191
$ws.Cells.Item(18, 8).Value = 5.745275844894876
$ws.Cells.Item(18, 9).Value = 7.20293239625024

# Row 19
$ws.Cells.Item(19, 1).Value = "model_2_7_8"
$ws.Cells.Item(19, 2).Value = 0.9808672478359243
$ws.Cells.Item(19, 3).Value = 0.9417513770558135
$ws.Cells.Item(19, 4).Value = 0.9498173348743632
$ws.Cells.Item(19, 5).Value = 0.9459883917421182
$ws.Cells.Item(19, 6).Value = 2.307415607393554
$ws.Cells.Item(19, 7).Value = 8.123750020751949
$ws.Cells.Item(19, 8).Value = 5.36584748485835
$ws.Cells.Item(19, 9).Value = 6.825912198673827

# Row 20
$ws.Cells.Item(20, 1).Value = "model_2_7_7"
$ws.Cells.Item(20, 2).Value = 0.981405544521755
$ws.Cells.Item(20, 3).Value = 0.9445463676330105
$ws.Cells.Item(20, 4).Value = 0.9535765892831496
$ws.Cells.Item(20, 5).Value = 0.9491180756258895
$ws.Cells.Item(20, 6).Value = 2.242496866814966
$ws.Cells.Item(20, 7).Value = 7.733941582168569
$ws.Cells.Item(20, 8).Value = 4.963884261824495
$ws.Cells.Item(20, 9).Value = 6.430387086771431

# Row 21
$ws.Cells.Item(21, 1).Value = "model_2_7_1"
$ws.Cells.Item(21, 2).Value = 0.9814543957031046
$ws.Cells.Item(21, 3).Value = 0.9607389167797603
$ws.Cells.Item(21, 4).Value = 0.9784321445976364
$ws.Cells.Item(21, 5).Value = 0.968474786414246
$ws.Cells.Item(21, 6).Value = 2.236605399799708
$ws.Cells.Item(21, 7).Value = 5.475618297977275
$ws.Cells.Item(21, 8).Value = 2.306171311842866
$ws.Cells.Item(21, 9).Value = 3.984112803184193

# Row 22
$ws.Cells.Item(22, 1).Value = "model_2_7_6"
$ws.Cells.Item(22, 2).Value = 0.9818751523578542
$ws.Cells.Item(22, 3).Value = 0.9474169748864287
$ws.Cells.Item(22, 4).Value = 0.9575281875877636
$ws.Cells.Item(22, 5).Value = 0.9523686304288548
$ws.Cells.Item(22, 6).Value = 2.185862021964776
$ws.Cells.Item(22, 7).Value = 7.33358712646113
$ws.Cells.Item(22, 8).Value = 4.541354414697046
$ws.Cells.Item(22, 9).Value = 6.019586475612432

# Row 23
$ws.Cells.Item(23, 1).Value = "model_2_7_2"
$ws.Cells.Item(23, 2).Value = 0.9821457312523677
$ws.Cells.Item(23, 3).Value = 0.9585242460154582
$ws.Cells.Item(23, 4).Value = 0.9743768539869582
$ws.Cells.Item(23, 5).Value = 0.9655662787918273
$ws.Cells.Item(23, 6).Value = 2.153230126726821
$ws.Cells.Item(23, 7).Value = 5.784491379572647
$ws.Cells.Item(23, 8).Value = 2.739788595205544
$ws.Cells.Item(23, 9).Value = 4.35168596569795

# Row 24
$ws.Cells.Item(24, 1).Value = "model_2_7_5"
$ws.Cells.Item(24, 2).Value = 0.9822389784257286
$ws.Cells.Item(24, 3).Value = 0.9503202203914954
$ws.Cells.Item(24, 4).Value = 0.9616414499483918
$ws.Cells.Item(24, 5).Value = 0.955702432938206
$ws.Cells.Item(24, 6).Value = 2.141984489856959
$ws.Cells.Item(24, 7).Value = 6.92868071769276
$ws.Cells.Item(24, 8).Value = 4.101538425707966
$ws.Cells.Item(24, 9).Value = 5.598265134690693

# Row 25
$ws.Cells.Item(25, 1).Value = "model_2_7_3"
$ws.Cells.Item(25, 2).Value = 0.9824441445936022
$ws.Cells.Item(25, 3).Value = 0.9559653201688975
$ws.Cells.Item(25, 4).Value = 0.9701462437594135
$ws.Cells.Item(25, 5).Value = 0.9623867590916715
$ws.Cells.Item(25, 6).Value = 2.11724138892716
$ws.Cells.Item(25, 7).Value = 6.141376621632698
$ws.Cells.Item(25, 8).Value = 3.192152159238121
$ws.Cells.Item(25, 9).Value = 4.753509259009163

# Row 26
$ws.Cells.Item(26, 1).Value = "model_2_7_4"
$ws.Cells.Item(26, 2).Value = 0.9824491129920742
$ws.Cells.Item(26, 3).Value = 0.953196036825369
$ws.Cells.Item(26, 4).Value = 0.9658696676609474
$ws.Cells.Item(26, 5).Value = 0.9590660940479229
$ws.Cells.Item(26, 6).Value = 2.116642198592196
$ws.Cells.Item(26, 7).Value = 6.527599754169469
$ws.Cells.Item(26, 8).Value = 3.649430684488002
$ws.Cells.Item(26, 9).Value = 5.173170305234832
